$d = $word.ActiveDocument
$rsq = [char]0x2019   # right single quotation mark (curly apostrophe) used in the original text

# --- Change 1: Update the "Course Schedule/Calendar" heading (not the TOC entry) ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Course Schedule/Calendar`r") {
        $rng = $p.Range
        $rng.End = $rng.End - 1   # exclude the paragraph mark
        $rng.Text = "Course Schedule/Calendar (AL)"
        break
    }
}

# --- Change 2: Merge the two body paragraphs following the heading into one, with new wording ---
# Remove the second paragraph entirely (its content is folded into the first, reworded).
$oldPara2 = "In addition to these functions, this menu will also have two additional redirects. One will download and upload a calendar file to Google Calendar, using Google" + $rsq + "s proprietary calendar API. The second will redirect you to the Campus Map section of the application."
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq ($oldPara2 + "`r")) {
        $p.Range.Delete()
        break
    }
}

# Replace the text of the first paragraph with the new merged wording.
$findRng = $d.Content
$found = $findRng.Find.Execute(
    "From the main menu, users can select the Course Schedule menu. This menu will present users with several options and displays. The first is the time of the course, followed by the section number and course name. This section will also display the start and end date of each course, as well as allow for users to search for classes to add to their current or future semesters.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
)
if ($found) {
    $findRng.Text = "From the main menu, users can select the Course Schedule menu. This menu will present users with several options and displays. The first is the time of the course, followed by the section number and course name. This section will also display the start and end date of each course, as well as allow for users to search for classes to add to their current or future semesters. Furthermore, there will be an option to display the current term courses on a calendar. This particular feature will leverage Google's Calendar API for its mobile friendly design and ease of use. In addition to these functions, this menu will also redirect users to the Campus Map section of the application."
}

# --- Change 3: Update the Campus Map directions sentence ---
$findRng2 = $d.Content
$found2 = $findRng2.Find.Execute(
    " will offer basic directions for the user to follow from their current location to the classroom.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
)
if ($found2) {
    $findRng2.Text = " will offer basic directions for the user to follow from their current location to rooms on campus ."
}
